$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2 (phone number) and G2 (date-like text) must stay literal text,
# otherwise Excel would auto-convert them into a number / date serial.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("G2").NumberFormat = "@"

$ws.Range("A2").Value = "Andrian"
$ws.Range("B2").Value = "Jock"
$ws.Range("C2").Value = "Dock"
$ws.Range("D2").Value = "13:00"
$ws.Range("E2").Value = "043764635645"
$ws.Range("F2").Value = "Tratarea cariei"
$ws.Range("G2").Value = "5/23/2024"
$ws.Range("H2").Value = "13:005/23/2024"
